# ZEV Jan R2-4: convert the roboticRNAPrep / RIBOSOMAL_BAND / SMALL_RNA_BANDS
# columns (H, I, K) on rows 2-20 from text ("no"/"Y") to real booleans
# (FALSE/TRUE), for consistency across the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 8).Value  = $false   # H: roboticRNAPrep -> FALSE
    $ws.Cells.Item($r, 9).Value  = $true    # I: RIBOSOMAL_BAND -> TRUE
    $ws.Cells.Item($r, 11).Value = $true    # K: SMALL_RNA_BANDS -> TRUE
}

# Leave the cursor where the author left it when saving.
$ws.Range("I24").Select()
